$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Responding to Uncertainty: The Importance of Covertness in Support for Retaliation to Cyber and Kinetic Attacks"
$ws.Range("F2").Value = "Open Access"

$ws.Range("B3").Value = "On 3D simultaneous attack against manoeuvring target with communication delays"
$ws.Range("F3").Value = "Open Access"

$ws.Range("B4").Value = "A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"

$ws.Range("B5").Value = "Fighting in Cyberspace: Internet Access and the Substitutability of Cyber and Military Operations"

$ws.Range("B6").Value = "Robust tracking strategy for nonlinear connected vehicle cyber-physical systems"

$ws.Range("B7").Value = "Digital Assays Part II: Digital Protein and Cell Assays"

$ws.Range("B8").Value = "On domains: Cyber and the practice of warfare"

$ws.Range("B9").Value = "Ontological security, cyber technology, and states’ responses"
$ws.Range("F9").Value = "Open Access"

$ws.Range("B10").Value = "Warring from the virtual to the real: Assessing the public’s threshold for war over cyber security"
$ws.Range("F10").Value = "Open Access"

$ws.Range("B11").Value = "A virtual necessity: Some modest steps toward greater cybersecurity"
